$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2022-08-10 20:58:01"

for ($row = 2; $row -le 73; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
